$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "502.47")
# are written verbatim instead of being parsed into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.718.68"
$ws.Range("E2").Value = "  -9.23%  "
$ws.Range("D3").Value = "3.156.56"
$ws.Range("E3").Value = "  -10.92%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "502.47"
$ws.Range("E5").Value = "  -9.98%  "
$ws.Range("D6").Value = "167.33"
$ws.Range("E6").Value = "  -15.98%  "
$ws.Range("D7").Value = "0.577"
$ws.Range("E7").Value = "  -11.05%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "3.153.65"
$ws.Range("E9").Value = "  -10.80%  "
$ws.Range("D10").Value = "0.581"
$ws.Range("E10").Value = "  -12.77%  "
$ws.Range("D11").Value = "53.27"
$ws.Range("E11").Value = "  -13.39%  "
$ws.Range("D12").Value = "0.127"
$ws.Range("E12").Value = "  -12.66%  "
$ws.Range("D13").Value = "0.0000245"
$ws.Range("E13").Value = "  -9.61%  "
$ws.Range("D14").Value = "8.61"
$ws.Range("E14").Value = "  -14.02%  "
$ws.Range("D15").Value = "3.684.96"
$ws.Range("E15").Value = "  -10.11%  "
$ws.Range("D16").Value = "3.169.29"
$ws.Range("E16").Value = "  -10.37%  "
$ws.Range("D17").Value = "61.823.60"
$ws.Range("E17").Value = "  -8.86%  "
$ws.Range("E18").Value = "  -10.92%  "
$ws.Range("D19").Value = "16.54"
$ws.Range("E19").Value = "  -10.56%  "
$ws.Range("D20").Value = "10.46"
$ws.Range("E20").Value = "  -12.56%  "
$ws.Range("D21").Value = "0.915"
$ws.Range("E21").Value = "  -11.71%  "
$ws.Range("D22").Value = "352.04"
$ws.Range("E22").Value = "  -11.99%  "
$ws.Range("D23").Value = "3.57"
$ws.Range("E23").Value = "  -11.24%  "
$ws.Range("D24").Value = "77.42"
$ws.Range("E24").Value = "  -10.01%  "
$ws.Range("D25").Value = "6.05"
$ws.Range("E25").Value = "  -1.92%  "
$ws.Range("D26").Value = "10.38"
$ws.Range("E26").Value = "  -12.58%  "
$ws.Range("D27").Value = "3.73"
$ws.Range("E27").Value = "  -3.51%  "
$ws.Range("D28").Value = "2.53"
$ws.Range("E28").Value = "  -11.27%  "
$ws.Range("D29").Value = "10.73"
$ws.Range("E29").Value = "  -14.04%  "
$ws.Range("D30").Value = "7.88"
$ws.Range("E30").Value = "  -12.09%  "
$ws.Range("D31").Value = "27.27"
$ws.Range("E31").Value = "  -13.20%  "
$ws.Range("D32").Value = "604.47"
$ws.Range("E32").Value = "  -16.76%  "
$ws.Range("D33").Value = "6.22"
$ws.Range("E33").Value = "  -12.82%  "
$ws.Range("D34").Value = "10.72"
$ws.Range("E34").Value = "  -9.33%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "0.0995"
$ws.Range("E36").Value = "  -11.18%  "
$ws.Range("D37").Value = "55.63"
$ws.Range("E37").Value = "  -13.67%  "
$ws.Range("D38").Value = "35.08"
$ws.Range("E38").Value = "  -9.72%  "
$ws.Range("D39").Value = "0.364"
$ws.Range("E39").Value = "  -8.05%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.40%  "
$ws.Range("D41").Value = "0.0₃0646"
$ws.Range("E41").Value = "  -6.62%  "
$ws.Range("D42").Value = "2.769.36"
$ws.Range("E42").Value = "  -10.44%  "
$ws.Range("D43").Value = "0.116"
$ws.Range("E43").Value = "  -12.30%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "2.56"
$ws.Range("E44").Value = "  -7.75%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "2.31"
$ws.Range("E45").Value = "  -8.29%  "
$ws.Range("D46").Value = "2.52"
$ws.Range("E46").Value = "  -16.85%  "
$ws.Range("D47").Value = "0.0372"
$ws.Range("E47").Value = "  -9.80%  "
$ws.Range("D48").Value = "2.88"
$ws.Range("E48").Value = "  -4.77%  "
$ws.Range("D49").Value = "2.64"
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").Value = "131.42"
$ws.Range("E50").Value = "  -5.98%  "
$ws.Range("D51").Value = "0.119"
$ws.Range("E51").Value = "  -12.26%  "

# Restore default cell style on column D (the text format above was only
# needed transiently to stop auto-numeric conversion on write).
$ws.Range("D2:D51").Style = "Normal"
